$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume/% change (E) columns for each coin row
$ws.Range("D2").Value = "29.449.57"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "1.876.97"
$ws.Range("E3").Value = "  +0.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7132"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.00"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3114"
$ws.Range("E8").Value = "  +0.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07719"
$ws.Range("E9").Value = "  -2.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.38"
$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08387"
$ws.Range("E11").Value = "  +1.55%  "

$ws.Range("D12").Value = "1.920.76"
$ws.Range("E12").Value = "  +2.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.257"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7181"
$ws.Range("E14").Value = "  -0.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.75"
$ws.Range("E15").Value = "  +1.04%  "

$ws.Range("D16").Value = "29.462.96"
$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008234"
$ws.Range("E17").Value = "  +5.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.990"
$ws.Range("E18").Value = "  +2.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.48"
$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9993"
$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.920"
$ws.Range("E23").Value = "  -0.98%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1621"
$ws.Range("E25").Value = "  +0.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.79"
$ws.Range("E26").Value = "  +0.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.039"
$ws.Range("E27").Value = "  +0.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.63"
$ws.Range("E28").Value = "  +2.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.508"
$ws.Range("E29").Value = "  +0.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.421"
$ws.Range("E30").Value = "  +0.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.298"
$ws.Range("E31").Value = "  -3.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.319"
$ws.Range("E32").Value = "  +5.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05233"
$ws.Range("E33").Value = "  +0.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.929"
$ws.Range("E34").Value = "  -0.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7767"
$ws.Range("E35").Value = "  +7.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.177"
$ws.Range("E36").Value = "  -0.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.680"
$ws.Range("E37").Value = "  +0.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01868"
$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.724"
$ws.Range("E39").Value = "  +1.15%  "

$ws.Range("D40").Value = "1.166.96"
$ws.Range("E40").Value = "  -0.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.422"
$ws.Range("E41").Value = "  +4.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.73"
$ws.Range("E42").Value = "  +1.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8922"
$ws.Range("E43").Value = "  -1.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.46"
$ws.Range("E44").Value = "  +2.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9995"
$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D46").Value = "2.029.54"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.799"
$ws.Range("E47").Value = "  +0.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5201"
$ws.Range("E48").Value = "  -1.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.418"
$ws.Range("E49").Value = "  +1.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4316"
$ws.Range("E50").Value = "  +0.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.080"
$ws.Range("E51").Value = "  +0.69%  "

# Rows 20 and 21: Avalanche and WrappedliquidstakedEther2.0 swap positions, with updated price/volume data
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.135.14"
$ws.Range("E20").Value = "  +0.99%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.25"
$ws.Range("E21").Value = "  +0.22%  "
